# Insert a new worksheet "2022_26" (20-week post-booster data) between
# "2022_06" and "2022_47", populated with KCOR summary data.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("2022_06")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "2022_26"

# ---- Header row (row 1): bold, centered, bordered -> copy style from an existing header cell
$headerVals = @("Dose_Combination", "YearOfBirth", "KCOR", "CI_Lower", "CI_Upper")
for ($c = 1; $c -le 5; $c++) {
    $ws.Cells.Item(1, $c).Value = $headerVals[$c - 1]
}
$wb.Worksheets.Item("2022_06").Range("A1:E1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)

# ---- Data rows 2-80: force text storage (matches source, which stores every value as a string)
# Column A never looks numeric; B:E contain year/KCOR/CI values that must stay text, not become numbers.
$ws.Range("B2:E80").NumberFormat = "@"

$rows = @(
    @("Reporting date: 2024-10-07", "", "", "", ""),
    @("1 vs 0", "", "", "", ""),
    @("", "ASMR (direct)", "0.9657", "0.915", "1.019"),
    @("", "All Ages", "0.8722", "0.828", "0.919"),
    @("", "1920", "0.4515", "0.358", "0.569"),
    @("", "1930", "1.2649", "1.138", "1.406"),
    @("", "1940", "1.2618", "1.139", "1.398"),
    @("", "1950", "0.5530", "0.489", "0.625"),
    @("", "1960", "1.0483", "0.899", "1.223"),
    @("", "1970", "1.4149", "1.163", "1.722"),
    @("", "1980", "2.9145", "2.274", "3.735"),
    @("", "1990", "0.3441", "0.219", "0.540"),
    @("", "2000", "0.8313", "0.427", "1.617"),
    @("", "", "", "", ""),
    @("2 vs 0", "", "", "", ""),
    @("", "ASMR (direct)", "1.0184", "1.003", "1.035"),
    @("", "All Ages", "1.0148", "0.999", "1.031"),
    @("", "1920", "0.9420", "0.879", "1.009"),
    @("", "1930", "1.0452", "1.011", "1.080"),
    @("", "1940", "1.0552", "1.025", "1.087"),
    @("", "1950", "0.9694", "0.938", "1.002"),
    @("", "1960", "0.8955", "0.856", "0.937"),
    @("", "1970", "1.1982", "1.128", "1.273"),
    @("", "1980", "1.1335", "1.031", "1.246"),
    @("", "1990", "0.8963", "0.780", "1.030"),
    @("", "2000", "1.0430", "0.865", "1.257"),
    @("", "", "", "", ""),
    @("2 vs 1", "", "", "", ""),
    @("", "ASMR (direct)", "1.0545", "0.999", "1.114"),
    @("", "All Ages", "1.1635", "1.104", "1.227"),
    @("", "1920", "2.0865", "1.651", "2.638"),
    @("", "1930", "0.8263", "0.743", "0.919"),
    @("", "1940", "0.8362", "0.754", "0.927"),
    @("", "1950", "1.7530", "1.549", "1.984"),
    @("", "1960", "0.8542", "0.732", "0.997"),
    @("", "1970", "0.8468", "0.695", "1.032"),
    @("", "1980", "0.3889", "0.302", "0.501"),
    @("", "1990", "2.6050", "1.651", "4.109"),
    @("", "2000", "0.7206", "0.370", "1.405"),
    @("", "", "", "", ""),
    @("3 vs 0", "", "", "", ""),
    @("", "ASMR (direct)", "0.9774", "0.966", "0.989"),
    @("", "All Ages", "0.9785", "0.967", "0.990"),
    @("", "1920", "0.7326", "0.699", "0.768"),
    @("", "1930", "0.9938", "0.971", "1.017"),
    @("", "1940", "1.0876", "1.064", "1.112"),
    @("", "1950", "0.8885", "0.866", "0.912"),
    @("", "1960", "0.9816", "0.944", "1.021"),
    @("", "1970", "1.0094", "0.954", "1.068"),
    @("", "1980", "1.0888", "0.989", "1.198"),
    @("", "1990", "0.8034", "0.689", "0.937"),
    @("", "2000", "1.5111", "1.177", "1.941"),
    @("", "", "", "", ""),
    @("3 vs 1", "", "", "", ""),
    @("", "ASMR (direct)", "1.0106", "0.958", "1.066"),
    @("", "All Ages", "1.1219", "1.065", "1.182"),
    @("", "1920", "1.6227", "1.291", "2.040"),
    @("", "1930", "0.7857", "0.708", "0.872"),
    @("", "1940", "0.8619", "0.779", "0.954"),
    @("", "1950", "1.6068", "1.422", "1.815"),
    @("", "1960", "0.9364", "0.803", "1.091"),
    @("", "1970", "0.7134", "0.586", "0.868"),
    @("", "1980", "0.3736", "0.290", "0.481"),
    @("", "1990", "2.3350", "1.473", "3.701"),
    @("", "2000", "0.8349", "0.420", "1.661"),
    @("", "", "", "", ""),
    @("3 vs 2", "", "", "", ""),
    @("", "ASMR (direct)", "0.9597", "0.947", "0.973"),
    @("", "All Ages", "0.9642", "0.952", "0.977"),
    @("", "1920", "0.7777", "0.732", "0.826"),
    @("", "1930", "0.9508", "0.925", "0.977"),
    @("", "1940", "1.0308", "1.006", "1.056"),
    @("", "1950", "0.9166", "0.891", "0.943"),
    @("", "1960", "1.0962", "1.050", "1.144"),
    @("", "1970", "0.8424", "0.793", "0.895"),
    @("", "1980", "0.9605", "0.863", "1.069"),
    @("", "1990", "0.8964", "0.759", "1.059"),
    @("", "2000", "1.4488", "1.121", "1.872"),
    @("", "", "", "", "")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $val = $rowData[$c - 1]
        if ($val -ne "") {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# Force the used-range / dimension to extend down to row 80 (last row is blank in the source).
# Paste the (default, unstyled) format of an already-touched, plain-style cell so this
# doesn't introduce a spurious new cell style.
$ws.Range("A2").Copy()
$ws.Range("A80").PasteSpecial(-4122)

$ws.Range("A1").Select()
